# Repull data, push all data, mean calculation
# Update dSF (column F) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -5
$ws.Range("F6").Value  = -1
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = -5
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -5
